$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    2 = @{ "D" = "245.76"; "E" = "1.00%"; "G" = "4" }
    3 = @{ "D" = "29.36"; "E" = "-2.49%"; "G" = "4" }
    4 = @{ "D" = "5.142"; "E" = "0.23%"; "G" = "4" }
    5 = @{ "D" = "0.05775"; "E" = "1.75%"; "G" = "4" }
    6 = @{ "D" = "6.608"; "E" = "1.38%"; "G" = "4" }
    7 = @{ "D" = "0.8578"; "E" = "3.01%"; "G" = "4" }
    8 = @{ "D" = "0.8592"; "E" = "-0.43%"; "G" = "4" }
    9 = @{ "D" = "0.1367"; "E" = "2.78%"; "G" = "4" }
    10 = @{ "D" = "0.07023"; "E" = "1.67%"; "G" = "4" }
    11 = @{ "D" = "0.02983"; "E" = "4.20%"; "G" = "4" }
    12 = @{ "D" = "0.09362"; "E" = "-0.16%"; "G" = "4" }
    13 = @{ "D" = "0.001524"; "E" = "0.99%"; "G" = "4" }
    14 = @{ "D" = "0.0006028"; "E" = "0.67%"; "G" = "4" }
    15 = @{ "D" = "0.006011"; "E" = "-0.15%"; "G" = "4" }
    16 = @{ "D" = "3.482"; "E" = "-1.12%"; "G" = "4" }
    17 = @{ "D" = "3.155"; "E" = "4.42%"; "G" = "4" }
    18 = @{ "D" = "2.172"; "E" = "1.91%"; "G" = "4" }
    19 = @{ "D" = "0.3202"; "E" = "1.68%"; "G" = "4" }
    20 = @{ "D" = "0.03316"; "E" = "1.21%"; "G" = "4" }
    21 = @{ "E" = "-0.71%"; "G" = "4" }
    22 = @{ "D" = "3.179"; "E" = "-12.43%"; "G" = "4" }
    23 = @{ "D" = "0.04145"; "E" = "-0.08%"; "G" = "4" }
    24 = @{ "D" = "0.1400"; "E" = "2.07%"; "G" = "4" }
    25 = @{ "E" = "1.39%"; "G" = "4" }
    26 = @{ "E" = "-7.14%"; "G" = "4" }
    27 = @{ "E" = "2.57%"; "G" = "4" }
    28 = @{ "D" = "0.0001449"; "E" = "3.38%"; "G" = "4" }
    29 = @{ "G" = "4" }
    30 = @{ "G" = "4" }
    31 = @{ "G" = "4" }
    32 = @{ "G" = "4" }
    33 = @{ "G" = "4" }
    34 = @{ "G" = "4" }
    35 = @{ "G" = "4" }
    36 = @{ "G" = "4" }
    37 = @{ "G" = "4" }
    38 = @{ "G" = "4" }
    39 = @{ "G" = "4" }
    40 = @{ "D" = "0.03732"; "E" = "0.58%"; "G" = "4" }
    41 = @{ "E" = "1.35%"; "G" = "4" }
    42 = @{ "E" = "5.66%"; "G" = "4" }
    43 = @{ "D" = "0.003505"; "E" = "-37.82%"; "G" = "4" }
    44 = @{ "D" = "0.008536"; "E" = "-12.86%"; "G" = "4" }
    45 = @{ "D" = "0.00005281"; "E" = "3.80%"; "G" = "4" }
    46 = @{ "E" = "0.18%"; "G" = "4" }
    47 = @{ "D" = "0.05798"; "E" = "-41.90%"; "G" = "4" }
    48 = @{ "E" = "-18.75%"; "G" = "4" }
    49 = @{ "E" = "0.18%"; "G" = "4" }
    50 = @{ "E" = "0.18%"; "G" = "4" }
    51 = @{ "G" = "4" }
}

foreach ($row in $updates.Keys) {
    $cols = $updates[$row]
    foreach ($col in $cols.Keys) {
        $cellRef = "$col$row"
        $cell = $ws.Range($cellRef)
        $cell.NumberFormat = "@"
        $cell.Value = $cols[$col]
    }
}
